$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.983.52"
$ws.Range("E2").Value = "  +0.26%  "

$ws.Range("D3").Value = "1.559.12"
$ws.Range("E3").Value = "  +0.53%  "

$ws.Range("E4").Value = "  -0.53%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "207.18"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.09%  "

$ws.Range("E6").Value = "  +1.02%  "

$ws.Range("E7").Value = "  -0.54%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "22.10"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.23%  "

$ws.Range("E9").Value = "  -0.36%  "

$ws.Range("E10").Value = "  +1.55%  "

$ws.Range("E11").Value = "  -0.30%  "

$ws.Range("D12").Value = "1.779.43"
$ws.Range("E12").Value = "  +0.46%  "

$ws.Range("D13").Value = "1.518.71"
$ws.Range("E13").Value = "  -2.05%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.77"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.40%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.522"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.45%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "61.98"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.36%  "

$ws.Range("D17").Value = "26.978.45"
$ws.Range("E17").Value = "  +0.27%  "

$ws.Range("B18").Value = "ShibaInu"
$ws.Range("C18").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D18").Value = "0.0₃0708"
$ws.Range("E18").Value = "  +2.84%  "

$ws.Range("B19").Value = "BitcoinCash"
$ws.Range("C19").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "218.07"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.42%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.33"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.58%  "

$ws.Range("E21").Value = "  -0.54%  "

$ws.Range("E22").Value = "  +2.11%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.23"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.38%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.92"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -3.02%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "153.43"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.77%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.65"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.04%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "15.05"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.23%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.105"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.34%  "

$ws.Range("E29").Value = "  -0.52%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0470"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.56%  "

$ws.Range("E31").Value = "  -0.11%  "

$ws.Range("E32").Value = "  +0.60%  "

$ws.Range("D33").Value = "1.427.74"
$ws.Range("E33").Value = "  +1.12%  "

$ws.Range("E34").Value = "  +4.27%  "

$ws.Range("E35").Value = "  +13.07%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.61"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +3.31%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.31"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.61%  "

$ws.Range("E38").Value = "  +0.32%  "

$ws.Range("E39").Value = "  +1.95%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.810"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.62%  "

$ws.Range("E41").Value = "  -0.52%  "

$ws.Range("E42").Value = "  +1.90%  "

$ws.Range("E43").Value = "  +2.27%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.00"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.96%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "64.80"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.00%  "

$ws.Range("E46").Value = "  -0.17%  "

$ws.Range("D47").Value = "1.693.42"
$ws.Range("E47").Value = "  +0.48%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "87.66"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.78%  "

$ws.Range("E49").Value = "  +1.26%  "

$ws.Range("D50").Value = "0.0₇0990"
$ws.Range("E50").Value = "  +2.58%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0961"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.73%  "
